$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A53").Value = "T come tigro"
$ws.Range("B53").Value = "MATTEO PILATI | Pinguini Trentini"
$ws.Range("C53").Value = "Leonardo Viola | Shark Attack"
$ws.Range("D53").Value = "Alessio Bragagna | FC Savignano"
$ws.Range("E53").Value = "Sebastiano Zoller | CGB Gamberoni"
$ws.Range("F53").Value = "Roberto Barozzi | Demobusters"
